$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fill in row 7 with the new work-log entry (11.03.2020)
$ws.Range("B7").Value = 43901
$ws.Range("C7").Value = 0.44444444444444442
$ws.Range("D7").Value = 0.51041666666666663

$ws.Range("F7").Value = "GitHub"
$ws.Range("G7").Value = "Bataille Navale"
$ws.Range("H7").Value = "Planifications des sprint"
$ws.Range("I7").Value = "Création des sprints pour chaques semaines"

# Update selection to match the saved view state
$ws.Range("B8").Select()

$wb.Save()
